$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.316.83"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "'1.879.36"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'237.31"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.4820"
$ws.Range("E7").Value = "  -2.33%  "
$ws.Range("D8").Value = "'0.2886"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("D9").Value = "'0.06586"
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D10").Value = "'1.880.12"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "'16.94"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "'0.07391"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'5.205"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'88.01"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "'0.6602"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "'30.272.12"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "'13.62"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'0.000007732"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "'5.468"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").Value = "'2.140.19"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'194.42"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").Value = "'9.443"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "'164.64"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'18.23"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("D28").Value = "'1.929"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "'1.443"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").Value = "'4.274"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").Value = "'0.09146"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "'4.051"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "'0.05066"
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("D34").Value = "'0.7413"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "'2.712"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'2.630"
$ws.Range("D39").Value = "'0.9159"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").Value = "'2.076"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "'106.24"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'5.882"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "'0.4324"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").Value = "'0.9991"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'7.653"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'0.1343"
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("D47").Value = "'1.580"
$ws.Range("E47").Value = "  +9.86%  "
$ws.Range("D48").Value = "'65.26"
$ws.Range("E48").Value = "  -10.64%  "
$ws.Range("D49").Value = "'8.898"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").Value = "'34.16"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "'0.05724"
$ws.Range("E51").Value = "  -2.62%  "
